$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview sheet) - mirrors de-de's Correspond Handoff Datetime
$wsOverview.Range("G2").Value = "2016-08-21 11:08:31"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-21 11:08:26"
$wsZhCn.Range("K2").Value = "2016-08-21 11:08:44"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-21 11:08:31"
$wsDeDe.Range("K2").Value = "2016-08-21 11:08:50"
